$d = $word.ActiveDocument
$dash = [char]0x2013

# ---------------------------------------------------------------------------
# 1) "Scores by [engnat]" paragraph: drop the spell-check markup around
#    "engnat" and merge the three runs into a single run, keeping the
#    leading tab as a distinct <w:tab/> run (not flattened into text).
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "`tScores by engnat " + $dash + " Carter`r") {
        $rng = $d.Range($p.Range.Start, $p.Range.End - 1)
        $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/><w:t>Scores by engnat ' + $dash + ' Carter</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $rng.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# 2) "Notable variables" paragraph: delete the trailing
#    " - Waiting on Nicholas, Carter" run entirely.
# ---------------------------------------------------------------------------
$search2 = " " + $dash + " Waiting on Nicholas, Carter"
$rng2 = $d.Content
$rng2.Find.Execute($search2, $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "S27, S26, " paragraph: replace with the new list of variables,
#    split across runs exactly as recorded in the target revision.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "S27, S26, `r") {
        $rng3 = $d.Range($p.Range.Start, $p.Range.End - 1)
        $xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">S5, S14, S16, </w:t></w:r><w:r><w:t>S2</w:t></w:r><w:r><w:t>6</w:t></w:r><w:r><w:t>, S2</w:t></w:r><w:r><w:t>7</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $rng3.InsertXML($xml3)
        break
    }
}
